$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "supprimereradmin" -> "supprimeradmin" in the "NOM DU USE CASE" cell (A1)
$cell = $ws.Range("A1")
$cell.Value = "NOM DU USE CASE : supprimeradmin"

# Restore the rich-text formatting: bold 12pt header followed by normal 11pt text.
$boldLen = "NOM DU USE CASE".Length
$cell.Characters(1, $boldLen).Font.Bold = $true
$cell.Characters(1, $boldLen).Font.Size = 12
$cell.Characters($boldLen + 1, $cell.Value.Length - $boldLen).Font.Bold = $false
$cell.Characters($boldLen + 1, $cell.Value.Length - $boldLen).Font.Size = 11

# Make sure the final selection is A1 (default), matching the diff which drops the
# stored <selection> element for the sheet view.
$ws.Range("A1").Select()
